# "Fruta / hortaliza, semanal" weekly refresh:
# the 19 data rows (2-20) on the active sheet get re-shuffled - each row's
# Fecha (D), Volumen (M), Precio minimo/maximo/promedio (N:P) and Precio $/Kg (S)
# are replaced with another row's values from that same weekly batch.
# Columns A,B,C,E-L,Q,R,T are untouched (constant across the sheet already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44294     # D2
$ws.Cells.Item(2, 13).Value = 25       # M2
$ws.Range("N2:P2").Value = 25000       # N2:P2
$ws.Cells.Item(2, 19).Value = 1250     # S2

$ws.Cells.Item(3, 4).Value = 44413     # D3
$ws.Cells.Item(3, 13).Value = 45       # M3
$ws.Range("N3:P3").Value = 20000       # N3:P3
$ws.Cells.Item(3, 19).Value = 1000     # S3

$ws.Cells.Item(4, 4).Value = 44291     # D4
$ws.Cells.Item(4, 13).Value = 70       # M4
$ws.Range("N4:P4").Value = 25000       # N4:P4
$ws.Cells.Item(4, 19).Value = 1250     # S4

$ws.Cells.Item(5, 4).Value = 44305     # D5
$ws.Cells.Item(5, 13).Value = 20       # M5
$ws.Range("N5:P5").Value = 22000       # N5:P5
$ws.Cells.Item(5, 19).Value = 1100     # S5

$ws.Cells.Item(6, 4).Value = 44292     # D6
$ws.Cells.Item(6, 13).Value = 30       # M6
$ws.Range("N6:P6").Value = 25000       # N6:P6
$ws.Cells.Item(6, 19).Value = 1250     # S6

$ws.Cells.Item(7, 4).Value = 44406     # D7
$ws.Cells.Item(7, 13).Value = 20       # M7
$ws.Range("N7:P7").Value = 20000       # N7:P7
$ws.Cells.Item(7, 19).Value = 1000     # S7

$ws.Cells.Item(8, 4).Value = 44385     # D8
$ws.Cells.Item(8, 13).Value = 36       # M8
$ws.Range("N8:P8").Value = 20000       # N8:P8
$ws.Cells.Item(8, 19).Value = 1000     # S8

$ws.Cells.Item(9, 4).Value = 44403     # D9
$ws.Cells.Item(9, 13).Value = 50       # M9
$ws.Range("N9:P9").Value = 20000       # N9:P9
$ws.Cells.Item(9, 19).Value = 1000     # S9

$ws.Cells.Item(10, 4).Value = 44307    # D10
$ws.Cells.Item(10, 13).Value = 30      # M10
$ws.Range("N10:P10").Value = 22000     # N10:P10
$ws.Cells.Item(10, 19).Value = 1100    # S10

$ws.Cells.Item(11, 4).Value = 44300    # D11
$ws.Cells.Item(11, 13).Value = 45      # M11
$ws.Range("N11:P11").Value = 22000     # N11:P11
$ws.Cells.Item(11, 19).Value = 1100    # S11

$ws.Cells.Item(12, 4).Value = 44298    # D12
$ws.Cells.Item(12, 13).Value = 65      # M12
$ws.Range("N12:P12").Value = 22000     # N12:P12
$ws.Cells.Item(12, 19).Value = 1100    # S12

$ws.Cells.Item(13, 4).Value = 44389    # D13
$ws.Cells.Item(13, 13).Value = 20      # M13
$ws.Range("N13:P13").Value = 20000     # N13:P13
$ws.Cells.Item(13, 19).Value = 1000    # S13

$ws.Cells.Item(14, 4).Value = 44301    # D14
$ws.Cells.Item(14, 13).Value = 38      # M14
$ws.Range("N14:P14").Value = 22000     # N14:P14
$ws.Cells.Item(14, 19).Value = 1100    # S14

$ws.Cells.Item(15, 4).Value = 44400    # D15
$ws.Cells.Item(15, 13).Value = 45      # M15
$ws.Range("N15:P15").Value = 20000     # N15:P15
$ws.Cells.Item(15, 19).Value = 1000    # S15

$ws.Cells.Item(16, 4).Value = 44382    # D16
$ws.Cells.Item(16, 13).Value = 24      # M16
$ws.Range("N16:P16").Value = 20000     # N16:P16
$ws.Cells.Item(16, 19).Value = 1000    # S16

$ws.Cells.Item(17, 4).Value = 44445    # D17
$ws.Cells.Item(17, 13).Value = 45      # M17
$ws.Range("N17:P17").Value = 20000     # N17:P17
$ws.Cells.Item(17, 19).Value = 1000    # S17

$ws.Cells.Item(18, 4).Value = 44377    # D18
$ws.Cells.Item(18, 13).Value = 25      # M18
$ws.Range("N18:P18").Value = 20000     # N18:P18
$ws.Cells.Item(18, 19).Value = 1000    # S18

$ws.Cells.Item(19, 4).Value = 44448    # D19
$ws.Cells.Item(19, 13).Value = 30      # M19
$ws.Range("N19:P19").Value = 22000     # N19:P19
$ws.Cells.Item(19, 19).Value = 1100    # S19

$ws.Cells.Item(20, 4).Value = 44376    # D20
$ws.Cells.Item(20, 13).Value = 38      # M20
$ws.Range("N20:P20").Value = 20000     # N20:P20
$ws.Cells.Item(20, 19).Value = 1000    # S20
